$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1515
$ws.Range("E2").Value = 76
$ws.Range("C3").Value = 152
$ws.Range("E5").Value = 89
$ws.Range("C6").Value = 34
$ws.Range("C7").Value = 33
$ws.Range("C8").Value = 34
$ws.Range("C9").Value = 34
$ws.Range("C10").Value = 34
